# Weekly update: a new "Ají" price record for Vega Monumental Concepción is
# inserted as row 65, pushing all subsequent rows (old 65..142) down by one
# (new rows 66..143). The worksheet's used range grows from A1:R142 to
# A1:R143 automatically once Excel inserts the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 65 - this shifts rows 65-142 down
# to 66-143, carrying their formatting (incl. the date-formatted column D).
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new record.
$ws.Range("A65").Value = 11
$ws.Range("B65").Value = "Vega Monumental Concepción"
$ws.Range("C65").Value = "Bíobío"
$ws.Range("D65").Value = "2022-11-09"
$ws.Range("E65").Value = 8
$ws.Range("F65").Value = 100112021
$ws.Range("G65").Value = "Ají"
$ws.Range("H65").Value = "Inferno"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 220
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = 15455
$ws.Range("N65").Value = "$/caja 10 kilos"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 1546
$ws.Range("Q65").Value = 10
$ws.Range("R65").Value = "Hortaliza"
